$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.005.86"
$ws.Range("E2").Value = "  -0.61%  "
$ws.Range("D3").Value = "'3.456.20"
$ws.Range("E3").Value = "  -0.88%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'579.19"
$ws.Range("E5").Value = "  -1.13%  "
$ws.Range("D6").Value = "'148.81"
$ws.Range("E6").Value = "  +0.54%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'0.480"
$ws.Range("E8").Value = "  +0.17%  "
$ws.Range("D9").Value = "'7.90"
$ws.Range("E9").Value = "  +2.58%  "
$ws.Range("E10").Value = "  -2.23%  "
$ws.Range("E11").Value = "  +2.47%  "
$ws.Range("D12").Value = "'4.051.17"
$ws.Range("E12").Value = "  -0.75%  "
$ws.Range("E13").Value = "  +2.45%  "
$ws.Range("D14").Value = "'28.47"
$ws.Range("E14").Value = "  -4.54%  "
$ws.Range("D15").Value = "'3.459.07"
$ws.Range("E15").Value = "  -0.76%  "
$ws.Range("E16").Value = "  -1.02%  "
$ws.Range("D17").Value = "'63.074.60"
$ws.Range("E17").Value = "  -0.49%  "
$ws.Range("D18").Value = "'6.49"
$ws.Range("E18").Value = "  +3.22%  "
$ws.Range("D19").Value = "'14.66"
$ws.Range("E19").Value = "  +2.14%  "
$ws.Range("D20").Value = "'9.21"
$ws.Range("E20").Value = "  -2.00%  "
$ws.Range("D21").Value = "'389.34"
$ws.Range("E21").Value = "  -0.96%  "
$ws.Range("E22").Value = "  -0.31%  "
$ws.Range("D23").Value = "'74.70"
$ws.Range("E23").Value = "  -0.70%  "
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("D25").Value = "'3.594.67"
$ws.Range("E25").Value = "  -0.99%  "
$ws.Range("E26").Value = "  -2.50%  "
$ws.Range("E27").Value = "  -1.31%  "
$ws.Range("E28").Value = "  -1.37%  "
$ws.Range("E29").Value = "  +0.15%  "
$ws.Range("D30").Value = "'8.07"
$ws.Range("E30").Value = "  -2.48%  "
$ws.Range("E31").Value = "  -1.81%  "
$ws.Range("E32").Value = "  +0.03%  "
$ws.Range("D33").Value = "'1.36"
$ws.Range("E33").Value = "  -5.01%  "
$ws.Range("D34").Value = "'23.36"
$ws.Range("E34").Value = "  -2.01%  "
$ws.Range("D35").Value = "'1.63"
$ws.Range("E35").Value = "  +3.61%  "
$ws.Range("E36").Value = "  +0.41%  "
$ws.Range("B37").Value = "EnergySwap"
$ws.Range("C37").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D37").Value = "'31.98"
$ws.Range("E37").Value = "  -1.49%  "
$ws.Range("B38").Value = "Aptos"
$ws.Range("C38").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D38").Value = "'7.06"
$ws.Range("E38").Value = "  -1.46%  "
$ws.Range("D39").Value = "'170.24"
$ws.Range("E39").Value = "  -0.54%  "
$ws.Range("D40").Value = "'3.495.16"
$ws.Range("E40").Value = "  -0.84%  "
$ws.Range("E41").Value = "  +1.30%  "
$ws.Range("E42").Value = "  -1.33%  "
$ws.Range("D43").Value = "'42.91"
$ws.Range("E43").Value = "  +1.10%  "
$ws.Range("E44").Value = "  -1.40%  "
$ws.Range("E45").Value = "  -2.97%  "
$ws.Range("D46").Value = "'1.18"
$ws.Range("E46").Value = "  -2.20%  "
$ws.Range("D47").Value = "'2.567.58"
$ws.Range("E47").Value = "  -1.92%  "
$ws.Range("D48").Value = "'2.29"
$ws.Range("E48").Value = "  -0.47%  "
$ws.Range("E49").Value = "  +1.89%  "
$ws.Range("D50").Value = "'22.69"
$ws.Range("E50").Value = "  -5.11%  "
$ws.Range("E51").Value = "  +0.07%  "
